# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet.
# Previously the last row (row 7, 2025-03-31) carried the "last row" date
# format (YYYY-MM-DD, no time). Now that a new last row (row 8,
# 2025-04-01) is appended, row 7 reverts to the regular date/time format
# used by all the other non-final rows, and row 8 takes on the
# "last row" date-only format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is no longer the final row -> give it the same number format as
# the other interior rows (A2:A6).
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 8: today's tallies.
$ws.Range("A8").Value = 45748
$ws.Range("B8").Value = 28
$ws.Range("C8").Value = 24
$ws.Range("D8").Value = 28

# Row 8 is now the final row -> apply the date-only "last row" format.
$ws.Range("A8").NumberFormat = "YYYY-MM-DD"
